$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "ThirdGitHubRepo"
$ws.Range("B10").Value = "Hfdlfsfdsf"
$ws.Range("C10").Value = "push"

$ws.Range("C10").Select()
